# Adds a new "Wind Turbine Onshore" process row to the SEC_Processes sheet
# and the corresponding row on the PP (Power Plants) sheet, then updates the
# active sheet / selections to match the author's final view state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. SEC_Processes: fill in row 10 with the new process definition
# ---------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("SEC_Processes")

$wsProc.Range("B10").Value = "ELE"
$wsProc.Range("D10").Value = "ELE_EX_WIND_TURBINE"
$wsProc.Range("E10").Value = "Wind Turbine Onshore"
$wsProc.Range("F10").Value = "PJ"
$wsProc.Range("G10").Value = "GWe"
$wsProc.Range("H10").Value = "DAYNITE"

# ---------------------------------------------------------------------
# 2. PP: fill in row 9 with the new power-plant technology, referencing
#    the freshly-created SEC_Processes row (row 10) and SEC_Comm row 9.
# ---------------------------------------------------------------------
$wsPP = $wb.Worksheets.Item("PP")

$wsPP.Range("B9").Formula = "=SEC_Processes!D10"
$wsPP.Range("C9").Formula = "=SEC_Processes!E10"
$wsPP.Range("D9").Formula = "=SEC_Comm!C9"

# E9 switches from a formula-driven cell to a literal value, and picks up
# the formatting already used by D8/E8 (fontId 5 / fillId 26, left-aligned).
$wsPP.Range("D8").Copy()
$wsPP.Range("E9").PasteSpecial(-4122)
$wsPP.Range("E9").Value = "ELEC_HV"

$wsPP.Range("F9").Value = 1.345
$wsPP.Range("G9").Value = 1

# H9 picks up the right-aligned numeric formatting already used by H8.
$wsPP.Range("H8").Copy()
$wsPP.Range("H9").PasteSpecial(-4122)
$wsPP.Range("H9").Value = 31.536

$wsPP.Range("I9").Value = 0.33
$wsPP.Range("J9").Value = 1

# ---------------------------------------------------------------------
# 3. View state: PP becomes the active/selected sheet (was MIN_IMP), with
#    selections updated on SEC_Processes and PP.
# ---------------------------------------------------------------------
[void]$wsProc.Range("H11").Select()

[void]$wsPP.Activate()
[void]$wsPP.Range("I10").Select()
